$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks; fresh ones will be added for G2:G7 below.
$ws.Range("A1:I1000").Hyperlinks.Delete()

# Row 2: Listed Companies / Circular-NSE
$ws.Range("A2").Value2 = "Listed Companies"
$ws.Range("B2").Value2 = "Circular-NSE"
$ws.Range("C2").Value2 = "'2025"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value2 = "October"
$ws.Range("E2").Value2 = "'2025-10-14"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value2 = "SEBI Circular on “Minimum information to be provided to the Audit Committee and Shareholders for approval of Related Party Transactions.”"
$ws.Range("H2").Value2 = "SEBI_Circular_on_Minimum_information_to_be_provided_to_the_Audit_Committee_and_Shareholders_for_appr.pdf"
$ws.Range("I2").Value2 = "/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/Listed Companies/Circular-NSE/2025/October/SEBI_Circular_on_Minimum_information_to_be_provided_to_the_Audit_Committee_and_Shareholders_for_appr.pdf"

# Row 3: Listed Companies / Circular-BSE
$ws.Range("A3").Value2 = "Listed Companies"
$ws.Range("B3").Value2 = "Circular-BSE"
$ws.Range("C3").Value2 = "'2025"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value2 = "October"
$ws.Range("E3").Value2 = "'2025-10-14"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value2 = "SEBI Circular on “Minimum information to be provided to the Audit Committee and Shareholders for approval of Related Party Transactions.”"
$ws.Range("H3").Value2 = "SEBI_Circular_on_Minimum_information_to_be_provided_to_the_Audit_Committee_and_Shareholders_for_appr.pdf"
$ws.Range("I3").Value2 = "/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/Listed Companies/Circular-BSE/2025/October/SEBI_Circular_on_Minimum_information_to_be_provided_to_the_Audit_Committee_and_Shareholders_for_appr.pdf"

# Row 4: AIF / Circulars
$ws.Range("A4").Value2 = "AIF"
$ws.Range("B4").Value2 = "Circulars"
$ws.Range("C4").Value2 = "'2025"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value2 = "October"
$ws.Range("E4").Value2 = "'2025-10-15"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value2 = "Relaxation in timeline for disclosure of allocation methodology by Angel Funds"
$ws.Range("H4").Value2 = "1760525216783.pdf"
$ws.Range("I4").Value2 = "/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/AIF/Circulars/2025/October/1760525216783.pdf"

# Row 5: SEBI / Circulars
$ws.Range("A5").Value2 = "SEBI"
$ws.Range("B5").Value2 = "Circulars"
$ws.Range("C5").Value2 = "'2025"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value2 = "October"
$ws.Range("E5").Value2 = "'2025-10-13"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value2 = "Minimum information to be provided to the Audit Committee and Shareholders for approval of Related Party Transactions"
$ws.Range("H5").Value2 = "1760356560260.pdf"
$ws.Range("I5").Value2 = "/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Circulars/2025/October/1760356560260.pdf"

# Row 6: SEBI / Consulatation Paper
$ws.Range("A6").Value2 = "SEBI"
$ws.Range("B6").Value2 = "Consulatation Paper"
$ws.Range("C6").Value2 = "'2025"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value2 = "October"
$ws.Range("E6").Value2 = "'2025-10-17"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value2 = "Consultation paper on proposed amendment to certain provisions of SEBI (LODR) Regulations, 2015 to facilitate transfer of securities transferred prior to April 1, 2019 and simplify the process of dematerialization of shares"
$ws.Range("H6").Value2 = "1760699641194.pdf"
$ws.Range("I6").Value2 = "/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Consulatation Paper/2025/October/1760699641194.pdf"

# Row 7: SEBI / Master circular
$ws.Range("A7").Value2 = "SEBI"
$ws.Range("B7").Value2 = "Master circular"
$ws.Range("C7").Value2 = "'2025"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value2 = "October"
$ws.Range("E7").Value2 = "'2025-10-15"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value2 = "Master Circular for issue and listing of Non-convertible Securities, Securitised Debt Instruments, Security Receipts, Municipal Debt Securities and Commercial Paper"
$ws.Range("H7").Value2 = "1760532257519.pdf"
$ws.Range("I7").Value2 = "/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Master circular/2025/October/1760532257519.pdf"

# Hyperlinks for column G (PDF_URL): add link, set display text, restore Hyperlink style
$ws.Hyperlinks.Add($ws.Range("G2"), "https://nsearchives.nseindia.com//web/circular/2025-10/NSE_Circular_14102025_20251014172300.pdf")
$ws.Range("G2").Value2 = "https://nsearchives.nseindia.com//web/circular/2025-10/NSE_Circular_14102025_20251014172300.pdf"
$ws.Range("G2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("G3"), "https://www.bseindia.com/markets/MarketInfo/DownloadAttach.aspx?id=20251014-34&attachedId=dbe355c0-4b37-4967-8429-b7b6c97dbeba")
$ws.Range("G3").Value2 = "https://www.bseindia.com/markets/MarketInfo/DownloadAttach.aspx?id=20251014-34&attachedId=dbe355c0-4b37-4967-8429-b7b6c97dbeba"
$ws.Range("G3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("G4"), "https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1760525216783.pdf")
$ws.Range("G4").Value2 = "https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1760525216783.pdf"
$ws.Range("G4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("G5"), "https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1760356560260.pdf")
$ws.Range("G5").Value2 = "https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1760356560260.pdf"
$ws.Range("G5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("G6"), "https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1760699641194.pdf")
$ws.Range("G6").Value2 = "https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1760699641194.pdf"
$ws.Range("G6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("G7"), "https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1760532257519.pdf")
$ws.Range("G7").Value2 = "https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1760532257519.pdf"
$ws.Range("G7").Style = "Hyperlink"

"done"